$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("D2")
$rng.NumberFormat = "@"
$rng.Value = '22.347.20'
$rng.Style = "Normal"
$ws.Range("E2").Value = '  +0.24%  '
$rng = $ws.Range("D3")
$rng.NumberFormat = "@"
$rng.Value = '1.565.59'
$rng.Style = "Normal"
$ws.Range("E3").Value = '  +0.28%  '
$ws.Range("E4").Value = '  +0.45%  '
$rng = $ws.Range("D5")
$rng.NumberFormat = "@"
$rng.Value = '1.010'
$rng.Style = "Normal"
$ws.Range("E5").Value = '  +0.82%  '
$rng = $ws.Range("D6")
$rng.NumberFormat = "@"
$rng.Value = '289.03'
$rng.Style = "Normal"
$ws.Range("E6").Value = '  -0.19%  '
$rng = $ws.Range("D7")
$rng.NumberFormat = "@"
$rng.Value = '0.3728'
$rng.Style = "Normal"
$ws.Range("E7").Value = '  +0.02%  '
$rng = $ws.Range("D8")
$rng.NumberFormat = "@"
$rng.Value = '49.11'
$rng.Style = "Normal"
$ws.Range("E8").Value = '  -0.39%  '
$rng = $ws.Range("D9")
$rng.NumberFormat = "@"
$rng.Value = '0.3358'
$rng.Style = "Normal"
$ws.Range("E9").Value = '  -1.41%  '
$rng = $ws.Range("D10")
$rng.NumberFormat = "@"
$rng.Value = '0.07404'
$rng.Style = "Normal"
$ws.Range("E10").Value = '  -3.01%  '
$rng = $ws.Range("D11")
$rng.NumberFormat = "@"
$rng.Value = '1.113'
$rng.Style = "Normal"
$ws.Range("E11").Value = '  -4.30%  '
$ws.Range("E12").Value = '  +0.47%  '
$rng = $ws.Range("D13")
$rng.NumberFormat = "@"
$rng.Value = '20.64'
$rng.Style = "Normal"
$ws.Range("E13").Value = '  -3.46%  '
$rng = $ws.Range("D14")
$rng.NumberFormat = "@"
$rng.Value = '5.833'
$rng.Style = "Normal"
$ws.Range("E14").Value = '  -3.14%  '
$rng = $ws.Range("D15")
$rng.NumberFormat = "@"
$rng.Value = '6.829'
$rng.Style = "Normal"
$ws.Range("E15").Value = '  -1.18%  '
$rng = $ws.Range("D16")
$rng.NumberFormat = "@"
$rng.Value = '1.578.66'
$rng.Style = "Normal"
$ws.Range("E16").Value = '  +1.16%  '
$ws.Range("E17").Value = '  -1.75%  '
$rng = $ws.Range("D18")
$rng.NumberFormat = "@"
$rng.Value = '88.84'
$rng.Style = "Normal"
$ws.Range("E18").Value = '  -0.99%  '
$rng = $ws.Range("D19")
$rng.NumberFormat = "@"
$rng.Value = '0.06676'
$rng.Style = "Normal"
$ws.Range("E19").Value = '  -0.68%  '
$rng = $ws.Range("D20")
$rng.NumberFormat = "@"
$rng.Value = '1.007'
$rng.Style = "Normal"
$ws.Range("E20").Value = '  +0.51%  '
$rng = $ws.Range("D21")
$rng.NumberFormat = "@"
$rng.Value = '6.106'
$rng.Style = "Normal"
$ws.Range("E21").Value = '  -2.00%  '
$rng = $ws.Range("D22")
$rng.NumberFormat = "@"
$rng.Value = '16.11'
$rng.Style = "Normal"
$ws.Range("E22").Value = '  -2.61%  '
$rng = $ws.Range("D23")
$rng.NumberFormat = "@"
$rng.Value = '11.77'
$rng.Style = "Normal"
$ws.Range("E23").Value = '  -1.56%  '
$rng = $ws.Range("D24")
$rng.NumberFormat = "@"
$rng.Value = '22.348.06'
$rng.Style = "Normal"
$ws.Range("E24").Value = '  +0.23%  '
$rng = $ws.Range("D25")
$rng.NumberFormat = "@"
$rng.Value = '2.370'
$rng.Style = "Normal"
$ws.Range("E25").Value = '  -1.26%  '
$rng = $ws.Range("D26")
$rng.NumberFormat = "@"
$rng.Value = '2.485'
$rng.Style = "Normal"
$ws.Range("E26").Value = '  -11.30%  '
$rng = $ws.Range("D27")
$rng.NumberFormat = "@"
$rng.Value = '19.75'
$rng.Style = "Normal"
$ws.Range("E27").Value = '  -2.01%  '
$rng = $ws.Range("D28")
$rng.NumberFormat = "@"
$rng.Value = '147.11'
$rng.Style = "Normal"
$ws.Range("E28").Value = '  +0.81%  '
$rng = $ws.Range("D29")
$rng.NumberFormat = "@"
$rng.Value = '4.994'
$rng.Style = "Normal"
$ws.Range("E29").Value = '  +0.33%  '
$rng = $ws.Range("D30")
$rng.NumberFormat = "@"
$rng.Value = '123.69'
$rng.Style = "Normal"
$ws.Range("E30").Value = '  -1.27%  '
$rng = $ws.Range("D31")
$rng.NumberFormat = "@"
$rng.Value = '1.749.01'
$rng.Style = "Normal"
$ws.Range("E31").Value = '  +0.78%  '
$ws.Range("E32").Value = '  -1.34%  '
$rng = $ws.Range("D33")
$rng.NumberFormat = "@"
$rng.Value = '0.9679'
$rng.Style = "Normal"
$ws.Range("E33").Value = '  -4.01%  '
$rng = $ws.Range("D34")
$rng.NumberFormat = "@"
$rng.Value = '5.840'
$rng.Style = "Normal"
$ws.Range("E34").Value = '  -5.18%  '
$rng = $ws.Range("D35")
$rng.NumberFormat = "@"
$rng.Value = '9.622'
$rng.Style = "Normal"
$ws.Range("E35").Value = '  -3.79%  '
$rng = $ws.Range("D36")
$rng.NumberFormat = "@"
$rng.Value = '0.08386'
$rng.Style = "Normal"
$ws.Range("E36").Value = '  -1.21%  '
$rng = $ws.Range("D37")
$rng.NumberFormat = "@"
$rng.Value = '1.384'
$rng.Style = "Normal"
$ws.Range("E37").Value = '  +5.95%  '
$rng = $ws.Range("D38")
$rng.NumberFormat = "@"
$rng.Value = '0.02442'
$rng.Style = "Normal"
$ws.Range("E38").Value = '  -3.97%  '
$rng = $ws.Range("D39")
$rng.NumberFormat = "@"
$rng.Value = '0.2244'
$rng.Style = "Normal"
$ws.Range("E39").Value = '  -2.76%  '
$rng = $ws.Range("D40")
$rng.NumberFormat = "@"
$rng.Value = '0.06326'
$rng.Style = "Normal"
$ws.Range("E40").Value = '  -0.87%  '
$rng = $ws.Range("D41")
$rng.NumberFormat = "@"
$rng.Value = '5.318'
$rng.Style = "Normal"
$ws.Range("E41").Value = '  -2.93%  '
$rng = $ws.Range("D42")
$rng.NumberFormat = "@"
$rng.Value = '0.6132'
$rng.Style = "Normal"
$ws.Range("E42").Value = '  -3.15%  '
$rng = $ws.Range("D43")
$rng.NumberFormat = "@"
$rng.Value = '10.88'
$rng.Style = "Normal"
$ws.Range("E43").Value = '  -6.73%  '
$ws.Range("E44").Value = '  +0.44%  '
$rng = $ws.Range("D45")
$rng.NumberFormat = "@"
$rng.Value = '13.82'
$rng.Style = "Normal"
$ws.Range("E45").Value = '  -2.42%  '
$rng = $ws.Range("D46")
$rng.NumberFormat = "@"
$rng.Value = '3.777'
$rng.Style = "Normal"
$ws.Range("E46").Value = '  +0.60%  '
$rng = $ws.Range("D47")
$rng.NumberFormat = "@"
$rng.Value = '0.5715'
$rng.Style = "Normal"
$ws.Range("E47").Value = '  -4.12%  '
$rng = $ws.Range("D48")
$rng.NumberFormat = "@"
$rng.Value = '2.019'
$rng.Style = "Normal"
$ws.Range("E48").Value = '  -3.13%  '
$rng = $ws.Range("D49")
$rng.NumberFormat = "@"
$rng.Value = '125.43'
$rng.Style = "Normal"
$ws.Range("E49").Value = '  +0.85%  '
$rng = $ws.Range("D50")
$rng.NumberFormat = "@"
$rng.Value = '1.221'
$rng.Style = "Normal"
$ws.Range("E50").Value = '  -3.42%  '
$rng = $ws.Range("D51")
$rng.NumberFormat = "@"
$rng.Value = '0.07286'
$rng.Style = "Normal"
$ws.Range("E51").Value = '  +0.19%  '
